$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 13 (shifts old rows 13-23 down to 14-24),
# matching the author re-inserting the "Docentes responsaveis" value row
# that sits between "Docentes responsaveis:" (row 12) and "Programa resumido:" (row 14).
$ws.Rows.Item(13).Insert()

# The inserted row copied column A formatting into A13; remove that stray cell
# so row 13 has no value in column A (matches target layout).
$ws.Range("A13").Clear()

# Give B13/C13 the same number/text formatting as the row below (style carried
# by columns B/C throughout the sheet) before writing the values into them.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("B13").Value = "2341641 - Maria da Rosa Capri"
$ws.Range("C13").Value = "2341641 - Maria da Rosa Capri"

# Fill in the real "Objetivos" (Objectives) text, replacing the placeholder
# teacher-name value that had erroneously been left in B10/C10.
$objetivos = "Apresentar aos alunos as bases teóricas e experimentais dos métodos instrumentais (quantitativos e qualitativos) de uso mais frequente na área química, incluindo o preparo de amostras e a criteriosa avaliação dos resultados analíticos. Ao final da disciplina, o aluno deve ser capaz de escolher e aplicar a metodologia mais adequada à solução dos problemas analíticos em geral, assim como interpretar resultados de análises químicas."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Row 14 = "Programa resumido:" -- replace placeholder "Semestral" with the real short syllabus.
$programaResumido = "Introdução à Análise Instrumental. Preparo de amostras. Métodos Espectroanalíticos: UV/Visível, Absorção Atômica, Emissão Atômica, Infravermelho. Métodos Eletroanalíticos: Potenciometria e Condutimetria. Métodos Cromatográficos: Cromatografia a Gás e Cromatografia Líquida de Alta Eficiência."
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Row 16 = "Programa:" -- fill in the full syllabus text (previously held a stray date value).
$programaCompleto = "1) Introdução à Análise Instrumental. Correlação entre métodos analíticos instrumentais e por via úmida. Preparo de amostras em meio sólido e em meios líquidos aquosos e não aquosos. Solubilização, digestão, fontes de energia aplicadas ao preparo, estabilização de amostras.`n2) Introdução aos Métodos Espectroanalíticos: Natureza da energia radiante. Espectro eletromagnético. Interação da radiação com a matéria. Absorção seletiva. Absortividade. Lei de Beer-Lambert. Curvas analíticas. `n3) Introdução à Espectrofotometria no UV/Visível. Instrumentação. Aplicações e interpretação de resultados. Determinações simultâneas. Parte Experimental.`n4) Introdução às Espectrometrias de Absorção e de Emissão Atômicas. Instrumentação. Interferências. Origem do espectro de emissão atômica. Fontes de atomização e de excitação. Calibração. Aplicações e interpretação de resultados. Parte Experimental.`n5) Introdução à Espectroscopia no Infravermelho. Instrumentação. Interpretação de espectros. Aplicações. Parte Experimental.`n6) Introdução aos Métodos Eletroanalíticos: Potenciometria e Condutimetria. Instrumentação.  Métodos diretos e indiretos. Aplicações e interpretação de resultados. Parte experimental.`n7) Introdução aos Métodos Cromatográficos. Conceitos básicos dos métodos de separação. Fases móvel e estacionária. Cromatografia planar em papel e em camada delgada. Cromatografia em coluna: cromatografia a gás e cromatografia líquida de alta eficiência. Instrumentação. Aplicações e interpretação de resultados. Parte Experimental."
$ws.Range("B16").Value = $programaCompleto
$ws.Range("C16").Value = $programaCompleto

# Row 19 = "Metodo:" -- fill in the evaluation-method text (previously held a stray teacher name).
$metodo = "A avaliação da disciplina será feita por meio de avaliações escritas individuais (provas) e avaliações de atividades em grupo (relatórios das aulas práticas e/ou trabalhos escritos e/ou apresentações de seminários)."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Row 20 = "Criterio:" -- fill in the grading-criteria text.
$criterio = "A Média Final (MF) será calculada pela média entre todas as avaliações realizadas durante o semestre, sendo o conjunto das avaliações individuais correspondentes a 75% da composição de MF e o conjunto das avaliações em grupo correspondentes a 25% da composição de MF. Será aprovado o aluno que obtiver MF maior ou igual a cinco e frequência mínima de 70% no semestre."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Row 21 = "Norma de recuperacao:" -- fill in the makeup-exam rules text.
$norma = "No período de Recuperação haverá horário previamente definido para resolução de dúvidas e será realizada uma avaliação escrita individual (Prova da Recuperação = PR), com conteúdo de todos os tópicos apresentados na disciplina durante o semestre.`nA Nota de Recuperação (NR) será dada pela média aritmética entre a Média do Semestre (MF) e a Prova da Recuperação (PR), sendo considerado aprovado o aluno que obtiver NR maior ou igual a cinco."
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# Row 22 = "Bibliografia:" -- fill in the bibliography text.
$bibliografia = "1) Skoog, D.A.; Holler, F.J. ; Nieman, T.A. Princípios de análise instrumental. 5. ed. Porto Alegre: Bookman,  2002.`n2) MENDHAM,J.; DENNEY, R.C.; BARNES, J.D. ; Thomas, M. Vogel: análise química quantitativa. 6. ed. Rio de Janeiro: Livros Técnicos e Científicos, 2002.`n3) OHLWEILER, O.A. Fundamentos de análise instrumental. Rio de Janeiro: Livros Técnicos e Científicos, 1981.`n4) KRUG, F.J. (org.) Métodos de preparo de amostras: fundamentos sobre métodos de preparo de amostras orgânicas e inorgânicas para análise elementar. 1. ed. Piracicaba: Edição do autor, 2008. `n5) COLLINS, C.H.; BRAGA, G.L.; BONATO, P.S. (Org.) Fundamentos de cromatografia. 1. ed. Campinas: Editora da UNICAMP, 2006.`nBibliografia complementar`n1) CHRISTIAN, G.D. Analytical chemistry. 4. ed. Nova York: John Wiley & Sons, 1986.`n2) DYER, J.R. Aplicação da espectroscopia de absorção aos compostos orgânicos. 1. Reimpressão. São Paulo: Edgard Blucher, 1977.`n3) SILVERSTEIN, R.M.; WEBSTER, F.X.; KIEMLE, D.J. Identificação espectrométrica de compostos orgânicos. 7. ed. Rio de Janeiro: Livros Técnicos e Científicos, 2007.`n4) WILLARD, H.H.; MERRITE, L.; DEAB, J. Instrumentação analítica. Lisboa: Fundação Calouste Gulbekian,  1989."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

